# Presentation added as PDF
#
# Adds two new slides to the end of the deck (positions 5 and 6), both
# using the "Title and Content" layout (the same CustomLayout already
# used by slides 2-4 -- ppLayout index 2 / slideLayout2.xml):
#   Slide 5 - "Related Patterns"   (Strategy Pattern / Memento Pattern)
#             plus the small corner logo picture reused from slide 2.
#   Slide 6 - "S.O.L.I.D principles" (five SOLID bullet lines)
#
# Helper: appends a new paragraph of text (with its own leading line
# break) to a TextRange and tags it with the deck's Danish locale, the
# same way the rest of the existing slides are tagged.
function Add-Paragraph {
    param($TextRange, [string]$Text)
    $added = $TextRange.InsertAfter("`r" + $Text)
    $added.LanguageID = "da-DK"
    return $added
}

# Helper: appends a new run of text to a TextRange *within the same
# paragraph* (no line break) and tags it with the deck's Danish locale.
function Add-Run {
    param($TextRange, [string]$Text)
    $added = $TextRange.InsertAfter($Text)
    $added.LanguageID = "da-DK"
    return $added
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 - "Related Patterns"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)

# Title placeholder
$title5 = $s5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "Related"
$title5.LanguageID = "da-DK"
Add-Run $title5 " Patterns" | Out-Null

# Content placeholder
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Strategy Pattern"
$body5.LanguageID = "da-DK"
Add-Paragraph $body5 "Memento Pattern" | Out-Null

# Slide-number placeholder (matches slides 1-4)
$s5.HeadersFooters.SlideNumber.Visible = $true
$s5.Shapes.Item(3).Name = "Pladsholder til slidenummer 3"

# Corner logo picture - reuse the one already embedded on slide 2 so
# the image media is shared rather than re-embedded.
$logoSrc = $p.Slides.Item(2).Shapes.Item(3)
$logoSrc.Copy()
$s5.Shapes.Paste() | Out-Null
$s5.Shapes.Item(4).Name = "Billede 4"

# ---------------------------------------------------------------------
# Slide 6 - "S.O.L.I.D principles"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)

# Title placeholder
$title6 = $s6.Shapes.Item(1).TextFrame.TextRange
$title6.Text = "S.O.L.I.D "
$title6.LanguageID = "da-DK"
Add-Run $title6 "principles" | Out-Null

# Content placeholder - five SOLID bullet lines
$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "S: SRP – Single Responsibility Principle"
$body6.LanguageID = "da-DK"
Add-Paragraph $body6 "O: OCP – Open-Closed Principle" | Out-Null
Add-Paragraph $body6 "L: LSP – Liskov Subsitution Principle" | Out-Null
Add-Paragraph $body6 "I: ISP – Interface Segregation Principle" | Out-Null
Add-Paragraph $body6 "D: DIP – Dependency Inversion Principle" | Out-Null

# Slide-number placeholder (matches slides 1-4)
$s6.HeadersFooters.SlideNumber.Visible = $true
$s6.Shapes.Item(3).Name = "Pladsholder til slidenummer 3"

Write-Host "Slides now: $($p.Slides.Count)"
